$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rows 15 & 16 (currently blank) become new TODO entries, styled like row 14.
$ws.Range("B14:D14").Copy()
$ws.Range("B15:D15").PasteSpecial(-4122)
$ws.Range("B14:D14").Copy()
$ws.Range("B16:D16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(15, 2).Value = "Entropy Bottle Neck"
$ws.Cells.Item(15, 4).Value = "Binary search of INCREASE ENTROPY"
$ws.Cells.Item(15, 3).Value = "TODO"

# 2. Status column: "OK" -> "DONE" for the existing status rows (5-10)
for ($r = 5; $r -le 10; $r++) {
    $ws.Cells.Item($r, 3).Value = "DONE"
}

# 3. Rows 12 & 13 move from the "TODO list" styling to the "status table"
#    styling (same look as rows 5-10), and their Status cell becomes "DONE".
$ws.Range("B5:D5").Copy()
$ws.Range("B12:D12").PasteSpecial(-4122)
$ws.Range("B5:D5").Copy()
$ws.Range("B13:D13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(12, 3).Value = "DONE"
$ws.Cells.Item(13, 3).Value = "DONE"

# 4. Finish off the second new TODO row.
$ws.Cells.Item(16, 2).Value = "PCA - Implement"
$ws.Cells.Item(16, 4).Value = "todo"
$ws.Cells.Item(16, 3).Value = "TODO"

# 5. Move the active selection like the author did while working.
$ws.Range("D17").Select()
